# Adds the "Coarse Grain" commodity data: extends rice/wheat/rra sheets
# with From State / To State / Values columns + data rows, and adds two
# new sheets ("coarse grain" and "coarse_grain") with their own data.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Font.Bold = $true
}

function Write-Row($ws, $row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "rice" — already had From/From State/To/To State/Commodity;
# add the "Values" column and six data rows.
# ---------------------------------------------------------------------
$wsRice = $wb.Worksheets.Item("rice")
Set-HeaderCell $wsRice 1 6 "Values"

$riceRows = @(
    @("BXC", "Haryana", "NNA",  "Bihar", "Rice", 0.5),
    @("BXC", "Haryana", "NNA",  "Bihar", "Rice", 0.5),
    @("BHT", "Haryana", "DMSJ", "Bihar", "Rice", 1),
    @("BHT", "Haryana", "NNA",  "Bihar", "Rice", 0.5),
    @("BHT", "Haryana", "DMSJ", "Bihar", "Rice", 1),
    @("BHT", "Haryana", "NNA",  "Bihar", "Rice", 0.5)
)
$r = 2
foreach ($row in $riceRows) {
    Write-Row $wsRice $r $row
    $r++
}

# ---------------------------------------------------------------------
# Sheet 2: "wheat" — only had From/To/Commodity; rebuild the header row
# with From State / To State / Values inserted, and add the data rows.
# ---------------------------------------------------------------------
$wsWheat = $wb.Worksheets.Item("wheat")
$wsWheat.Cells.Clear()

$headers = @("From", "From State", "To", "To State", "Commodity", "Values")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $txt = $headers[$i]
    Set-HeaderCell $wsWheat 1 $col $txt
}

$wheatRows = @(
    @("NMH", "MP",     "JPTN", "Andhra Pradesh", "Wheat", 1),
    @("KSA", "Punjab", "FCIG", "Andhra Pradesh", "Wheat", 2)
)
$r = 2
foreach ($row in $wheatRows) {
    Write-Row $wsWheat $r $row
    $r++
}

# ---------------------------------------------------------------------
# Sheet 3: "rra" — only had From/To/Commodity; rebuild the header row
# the same way, and replace the old sample row with the new data.
# ---------------------------------------------------------------------
$wsRra = $wb.Worksheets.Item("rra")
$wsRra.Cells.Clear()

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $txt = $headers[$i]
    Set-HeaderCell $wsRra 1 $col $txt
}

Write-Row $wsRra 2 @("BHT", "Haryana", "NNA", "Bihar", "RRA", 2)

# ---------------------------------------------------------------------
# New sheets: "coarse grain" and "coarse_grain", appended after "rra",
# each with the same 6-column header and one data row.
# ---------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIdx)
$wsCoarse1 = $wb.Worksheets.Add($null, $afterSheet)
$wsCoarse1.Name = "coarse grain"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $txt = $headers[$i]
    Set-HeaderCell $wsCoarse1 1 $col $txt
}
Write-Row $wsCoarse1 2 @("RJY", "Haryana", "CTO", "Bihar", "Coarse Grain", 1)

$lastIdx2 = $wb.Worksheets.Count
$afterSheet2 = $wb.Worksheets.Item($lastIdx2)
$wsCoarse2 = $wb.Worksheets.Add($null, $afterSheet2)
$wsCoarse2.Name = "coarse_grain"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $txt = $headers[$i]
    Set-HeaderCell $wsCoarse2 1 $col $txt
}
Write-Row $wsCoarse2 2 @("RJY", "Haryana", "CTO", "Bihar", "Coarse Grain", 1)
